$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 133, shifting the existing rows 133:146 down to 134:147.
# Excel's default row insert copies formatting from the row above, matching
# the source workbook (date column D keeps its date style).
$ws.Rows(133).Insert()

# Populate the new row 133 with the new weekly record. Columns that are
# identical to the neighboring rows (A, B, C, E, F, G, H, I, N, O, Q, R)
# are copied from row 134 (the row that used to be row 133).
$ws.Range("A133").Value = $ws.Range("A134").Value()
$ws.Range("B133").Value = $ws.Range("B134").Value()
$ws.Range("C133").Value = $ws.Range("C134").Value()
$ws.Range("D133").Value = 45166
$ws.Range("E133").Value = $ws.Range("E134").Value()
$ws.Range("F133").Value = $ws.Range("F134").Value()
$ws.Range("G133").Value = $ws.Range("G134").Value()
$ws.Range("H133").Value = $ws.Range("H134").Value()
$ws.Range("I133").Value = $ws.Range("I134").Value()
$ws.Range("J133").Value = 100
$ws.Range("K133").Value = 14000
$ws.Range("L133").Value = 14000
$ws.Range("M133").Value = 14000
$ws.Range("N133").Value = $ws.Range("N134").Value()
$ws.Range("O133").Value = $ws.Range("O134").Value()
$ws.Range("P133").Value = 1077
$ws.Range("Q133").Value = $ws.Range("Q134").Value()
$ws.Range("R133").Value = $ws.Range("R134").Value()
